$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# 1. Row 16 (CHE2): breakdown changes from "Switzerland" to "Switzerland, Zurich"
$ws.Range("B16").Value = "Switzerland, Zurich"

# 2. Insert a new row at 17 for SWE1 / Sweden
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "SWE1"
$ws.Range("B17").Value = "Sweden"
$ws.Range("C17").Value = "ageband"
$ws.Range("D17").Value = "data/derived/SWE1/SWE1_agebands.RDS"
$ws.Range("E17").Value = "yes"

# Old row 17 (KEN1/Kenya) is now row 18 - no changes needed there.

# 3. Insert a new row at 19 for LA_CA1 / USA, Los Angeles
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "LA_CA1"
$ws.Range("B19").Value = "USA, Los Angeles"
$ws.Range("C19").Value = "ageband"
$ws.Range("D19").Value = "data/derived/USA/LA_CA1_agebands.RDS"
$ws.Range("E19").Value = "yes"

# Old row 18 (NYS1 ageband) is now row 20 - update location text & stratified flag
$ws.Range("B20").Value = "USA, New York State"
$ws.Range("E20").Value = "no"

# Old row 19 (NYS1 region) is now row 21 - update location text & stratified flag
$ws.Range("B21").Value = "USA, New York State"
$ws.Range("E21").Value = "no"

# Old rows 20-23 are now rows 22-25 (values unchanged).
# Rows 24 and 25 (old rows 22 and 23) lose their explicit style on columns A:D.
$ws.Range("A24:D24").Style = "Normal"
$ws.Range("A25:D25").Style = "Normal"

# Update the view to match the saved state (active selection ends on the last edited cell).
$ws.Activate()
$ws.Range("E25").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
